$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "68.832.83"
$ws.Range("E2").Value2 = "  -0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.456.21"
$ws.Range("E3").Value2 = "  -1.30%  "
$ws.Range("E4").Value2 = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "558.27"
$ws.Range("E5").Value2 = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "162.54"
$ws.Range("E6").Value2 = "  -1.86%  "
$ws.Range("E8").Value2 = "  -1.21%  "
$ws.Range("E9").Value2 = "  -3.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "2.288.73"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.335"
$ws.Range("E12").Value2 = "  -3.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.82"
$ws.Range("E13").Value2 = "  -0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "2.905.81"
$ws.Range("E14").Value2 = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "68.620.99"
$ws.Range("E15").Value2 = "  -0.79%  "
$ws.Range("E16").Value2 = "  -2.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "23.59"
$ws.Range("E17").Value2 = "  -1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.448.06"
$ws.Range("E18").Value2 = "  -1.92%  "
$ws.Range("E19").Value2 = "  -3.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "341.23"
$ws.Range("E20").Value2 = "  -3.24%  "
$ws.Range("E21").Value2 = "  -3.91%  "
$ws.Range("E22").Value2 = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.95"
$ws.Range("E23").Value2 = "  +2.17%  "
$ws.Range("E24").Value2 = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "66.89"
$ws.Range("E25").Value2 = "  -3.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "3.70"
$ws.Range("E26").Value2 = "  -2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.580.87"
$ws.Range("E27").Value2 = "  -1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.997"
$ws.Range("E28").Value2 = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "8.20"
$ws.Range("E29").Value2 = "  -4.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.0₃0817"
$ws.Range("E30").Value2 = "  -5.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "7.18"
$ws.Range("E31").Value2 = "  -4.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "435.81"
$ws.Range("E32").Value2 = "  -0.30%  "
$ws.Range("E33").Value2 = "  -0.07%  "
$ws.Range("E34").Value2 = "  -3.37%  "
$ws.Range("E35").Value2 = "  -5.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "157.15"
$ws.Range("E36").Value2 = "  +1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "19.04"
$ws.Range("E37").Value2 = "  -0.03%  "
$ws.Range("E38").Value2 = "  +0.11%  "
$ws.Range("E39").Value2 = "  -3.69%  "
$ws.Range("E40").Value2 = "  -1.54%  "
$ws.Range("E41").Value2 = "  -2.55%  "
$ws.Range("E42").Value2 = "  -3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "37.44"
$ws.Range("E43").Value2 = "  -1.08%  "
$ws.Range("E44").Value2 = "  -5.45%  "
$ws.Range("E45").Value2 = "  +3.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.08"
$ws.Range("E46").Value2 = "  -4.28%  "
$ws.Range("E47").Value2 = "  -3.80%  "
$ws.Range("E48").Value2 = "  -2.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0719"
$ws.Range("E49").Value2 = "  -0.54%  "
$ws.Range("E50").Value2 = "  -4.14%  "
$ws.Range("E51").Value2 = "  -2.41%  "
